$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (old A -> B, old B -> C) and label it "Branch".
$ws.Columns("A").Insert()
$ws.Range("A1").Value = "Branch"

# Number the task groups in the new column.
$ws.Range("A2").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A2:A3").Merge()

# Mark the "Construction de l'arbre de codage" task as done.
$ws.Range("C4").Value = "Fait"

# Center the new Branch column.
$ws.Range("A1:A7").HorizontalAlignment = -4108
$ws.Range("A1:A7").VerticalAlignment = -4108

# Recreate the conditional formatting that used to target column B on column C.
$ws.Range("B2").FormatConditions.Delete()
$ws.Range("B3").FormatConditions.Delete()
$ws.Range("B1:B1048576").FormatConditions.Delete()

$fc = $ws.Range("C2").FormatConditions.Add(1, 3, '"En cours"')
$fc.Font.Color = 22428
$fc.Interior.Color = 10284031

$fc = $ws.Range("C2").FormatConditions.Add(1, 3, '"Fait"')
$fc.Interior.Color = 5287936

$fc = $ws.Range("C3").FormatConditions.Add(1, 3, '"En cours"')
$fc.Font.Color = 22428
$fc.Interior.Color = 10284031

$fc = $ws.Range("C1:C1048576").FormatConditions.Add(1, 3, '"Non fait"')
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615

$fc = $ws.Range("C1:C1048576").FormatConditions.Add(1, 3, '"Fait"')
$fc.Font.Color = 24832
$fc.Interior.Color = 13561798

$fc = $ws.Range("C1:C1048576").FormatConditions.Add(1, 3, '"En cours"')
$fc.Font.Color = 22428
$fc.Interior.Color = 10284031

$ws.Range("B15").Select()
